$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 423, shifting existing rows 423-503 down to rows 424-504
# (Excel default: new row inherits formatting from the row above it).
$ws.Rows.Item(423).Insert()

# Populate the newly inserted row 423 with the new weekly price record for Cebolla.
$ws.Cells.Item(423, 1).Value = 5
$ws.Cells.Item(423, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(423, 3).Value = "Maule"
$ws.Cells.Item(423, 4).Value = 44637
$ws.Cells.Item(423, 5).Value = 7
$ws.Cells.Item(423, 6).Value = 100112004
$ws.Cells.Item(423, 7).Value = "Cebolla"
$ws.Cells.Item(423, 8).Value = "Sin especificar"
$ws.Cells.Item(423, 9).Value = "1a (cosecha)"
$ws.Cells.Item(423, 10).Value = 2500
$ws.Cells.Item(423, 11).Value = 4500
$ws.Cells.Item(423, 12).Value = 4500
$ws.Cells.Item(423, 13).Value = 4500
$ws.Cells.Item(423, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(423, 15).Value = 'Región del Maule'
$ws.Cells.Item(423, 16).Value = 180
$ws.Cells.Item(423, 17).Value = 25
$ws.Cells.Item(423, 18).Value = "Hortaliza"
